# Append the 3-year return-period rows (lead time 3d, forecast dates
# 2025-10-30 and 2025-10-31) to the Panay / Dao Bridge flood trigger
# analysis sheet, mirroring the existing row layout (rows 2-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("G5369")

$newRows = @(
    @{
        Row = 31
        A = "Philippines"
        B = "philippines"
        C = "Panay River Basin"
        D = "panay"
        E = "Dao Bridge"
        F = "G5369"
        G = "primary"
        H = "2025-10-30"
        I = 3
        J = 11.4249999999999
        K = 122.7249999999997
        L = 5
        M = 864.7614412809821
        N = "LOW"
        O = 603.6038567117938
        P = 864.7614412809821
        Q = 50
        R = 0
        S = 0
        T = 187.921875
        U = 188.5568695068359
        V = 141.125
        W = 271.484375
        X = 169.521484375
        Y = 207.759765625
        Z = $false
        AA = -78.26893452584693
    },
    @{
        Row = 32
        A = "Philippines"
        B = "philippines"
        C = "Panay River Basin"
        D = "panay"
        E = "Dao Bridge"
        F = "G5369"
        G = "primary"
        H = "2025-10-31"
        I = 3
        J = 11.4249999999999
        K = 122.7249999999997
        L = 5
        M = 864.7614412809821
        N = "LOW"
        O = 603.6038567117938
        P = 864.7614412809821
        Q = 50
        R = 0
        S = 0
        T = 154.6328125
        U = 157.0551605224609
        V = 124.7890625
        W = 210.2890625
        X = 144.919921875
        Y = 163.240234375
        Z = $false
        AA = -82.1184427151446
    }
)

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

foreach ($rowData in $newRows) {
    $r = $rowData["Row"]
    foreach ($col in $colOrder) {
        $cell = $ws.Range($col + $r)
        if ($col -eq "H") {
            # forecast_date is stored as plain text ("YYYY-MM-DD"), not a
            # real date serial - force text formatting before assigning so
            # Excel doesn't auto-convert the string into a date value.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$col]
    }
}
